$d = $word.ActiveDocument

# The three inline pictures (Pearson logo in footer1/footer2, BTec logo in
# header1) have their wp:docPr/pic:cNvPr "name" attribute swapped:
#   footer1 / footer2 Pearson logo : image1.png -> image2.png
#   header1 BTec logo              : image2.jpg -> image1.jpg
#
# InlineShape has no scriptable "Name" setter that touches both the
# wp:docPr AND pic:cNvPr name attributes at once, so we round-trip the
# whole package through WordOpenXML and patch the name="..." attributes
# directly - this keeps both copies (docPr + cNvPr) in sync, exactly like
# the authored diff.

$xml = $d.WordOpenXML

$xml = $xml.Replace('name="image1.png"', 'name="image2.png"')
$xml = $xml.Replace('name="image2.jpg"', 'name="image1.jpg"')

$d.WordOpenXML = $xml

Write-Host "renamed inline image parts"
